# Added a new way to calculate the field using numpy so the calculations
# are roughly 45 times faster than the previous iteration.
#
# This negates every value in column K (rows 9:96) and adds X / Y header
# labels (with a yellow fill + centered alignment) in J8:K8. Column L
# already holds the formula "=K*-1" so it recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Negate every value in K9:K96 (L9:L96 = K*-1 recalculates automatically).
for ($r = 9; $r -le 96; $r++) {
    $cell = $ws.Cells.Item($r, 11)
    $cur = $cell.Value()
    $cell.Value = -$cur
}

# Add the new "X" / "Y" headers in J8:K8.
$ws.Range("J8").Value = "X"
$ws.Range("K8").Value = "Y"

# Style them: centered alignment + solid yellow fill.
$hdr = $ws.Range("J8:K8")
$hdr.HorizontalAlignment = -4108
$hdr.Interior.Color = 65535

# Move the active selection from N12 to N8.
[void]$ws.Range("N8").Select()
